$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 98: 'The Dotted Line' / 'Enchanted Durium Ink'
$ws.Range("H98").Value = 2354.3635
$ws.Range("I98").Value = 2589.3103
$ws.Range("J98").Value = 651
$ws.Range("K98").Value = 2589.3103
$ws.Range("L98").Value = 651
$ws.Range("M98").Value = -1091.3103
$ws.Range("N98").Value = -3647

# ALC row 122: 'Wishful Inking' / 'Enchanted High Durium Ink'
$ws.Range("H122").Value = 2354.3635
$ws.Range("I122").Value = 2589.3103
$ws.Range("J122").Value = 651
$ws.Range("K122").Value = 7767.9309
$ws.Range("L122").Value = 1953
$ws.Range("M122").Value = -5317.9309
$ws.Range("N122").Value = -6853

# ALC row 125: 'Body over Mind' / 'Grade 5 Dexterity Alkahest'
$ws.Range("H125").Value = 1909.6666
$ws.Range("J125").Value = 2155.5
$ws.Range("L125").Value = 19399.5
$ws.Range("N125").Value = -24319.5

# ALC row 137: 'Cutting Edge of Culinary Quality' / 'Magnesia Whetstone'
$ws.Range("H137").Value = 23257504
$ws.Range("I137").Value = 1208
$ws.Range("K137").Value = 3624
$ws.Range("M137").Value = -1074

# ALC row 138: 'All-night Crafting' / 'Cunning Craftsman''s Tisane'
$ws.Range("H138").Value = 2396.738
$ws.Range("I138").Value = 2355.84
$ws.Range("K138").Value = 7067.52
$ws.Range("M138").Value = -1927.52

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2: 'Ain''t Got No Ingots' / 'Bronze Ingot'
$ws.Range("H2").Value = 3268662.8
$ws.Range("I2").Value = 680
$ws.Range("J2").Value = 7353641
$ws.Range("K2").Value = 680
$ws.Range("L2").Value = 7353641
$ws.Range("M2").Value = -567
$ws.Range("N2").Value = -7353867

# ARM row 116: 'No Scope' / 'Titanbronze Ingot'
$ws.Range("H116").Value = 3268662.8
$ws.Range("I116").Value = 680
$ws.Range("J116").Value = 7353641
$ws.Range("K116").Value = 680
$ws.Range("L116").Value = 7353641
$ws.Range("M116").Value = 1614
$ws.Range("N116").Value = -7358229

# ARM row 122: 'Haste for High Durium' / 'High Durium Nugget'
$ws.Range("H122").Value = 1244.3182
$ws.Range("I122").Value = 1276.3158
$ws.Range("J122").Value = 1041.6666
$ws.Range("K122").Value = 3828.9474
$ws.Range("L122").Value = 3124.9998
$ws.Range("M122").Value = -1378.9474
$ws.Range("N122").Value = -8024.9998

# ARM row 128: 'Heading toward Bankruptcy' / 'Manganese Helm of the Falling Dragon'
$ws.Range("H128").Value = 46250
$ws.Range("J128").Value = 46250
$ws.Range("L128").Value = 46250
$ws.Range("N128").Value = -56210

# ARM row 132: 'Don''t Bore Me, Ore Me' / 'Mountain Chromite Ingot'
$ws.Range("H132").Value = 2681.4119
$ws.Range("I132").Value = 2614.5
$ws.Range("J132").Value = 2842
$ws.Range("K132").Value = 7843.5
$ws.Range("L132").Value = 8526
$ws.Range("M132").Value = -5313.5
$ws.Range("N132").Value = -13586

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3: 'Hells Bells' / 'Bronze Ingot'
$ws.Range("H3").Value = 3268662.8
$ws.Range("I3").Value = 680
$ws.Range("J3").Value = 7353641
$ws.Range("K3").Value = 680
$ws.Range("L3").Value = 7353641
$ws.Range("M3").Value = -566
$ws.Range("N3").Value = -7353869

# BSM row 98: 'Killer Cutlery' / 'Doman Iron Culinary Knife'
$ws.Range("H98").Value = 59390.5
$ws.Range("J98").Value = 59390.5
$ws.Range("L98").Value = 59390.5
$ws.Range("N98").Value = -65380.5

# BSM row 107: 'The Gold Experience' / 'Deepgold Nugget'
$ws.Range("H107").Value = 2357.1428
$ws.Range("I107").Value = 2225
$ws.Range("J107").Value = 2533.3333
$ws.Range("K107").Value = 2225
$ws.Range("L107").Value = 2533.3333
$ws.Range("M107").Value = -305
$ws.Range("N107").Value = -6373.3333

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16: 'Raise the Roof' / 'Ash Lumber'
$ws.Range("H16").Value = 2820
$ws.Range("I16").Value = 3250
$ws.Range("J16").Value = 1702
$ws.Range("K16").Value = 3250
$ws.Range("L16").Value = 1702
$ws.Range("M16").Value = -2963
$ws.Range("N16").Value = -2276

# CRP row 31: 'Wall Not Found' / 'Walnut Lumber'
$ws.Range("H31").Value = 1236.3939
$ws.Range("I31").Value = 1161.138
$ws.Range("J31").Value = 1782
$ws.Range("K31").Value = 1161.138
$ws.Range("L31").Value = 1782
$ws.Range("M31").Value = -866.1379999999999
$ws.Range("N31").Value = -2372

# CRP row 34: 'Armoires of the Rich and Famous' / 'Walnut Lumber'
$ws.Range("H34").Value = 1236.3939
$ws.Range("I34").Value = 1161.138
$ws.Range("J34").Value = 1782
$ws.Range("K34").Value = 1161.138
$ws.Range("L34").Value = 1782
$ws.Range("M34").Value = -959.1379999999999
$ws.Range("N34").Value = -2186

# CRP row 105: 'Zelkova, My Love' / 'Zelkova Lumber'
$ws.Range("H105").Value = 707.8570999999999
$ws.Range("I105").Value = 589.25
$ws.Range("J105").Value = 866
$ws.Range("K105").Value = 589.25
$ws.Range("L105").Value = 866
$ws.Range("M105").Value = 1157.75
$ws.Range("N105").Value = -4360

# CRP row 113: 'Patient Patients' / 'White Ash Lumber'
$ws.Range("H113").Value = 2820
$ws.Range("I113").Value = 3250
$ws.Range("J113").Value = 1702
$ws.Range("K113").Value = 3250
$ws.Range("L113").Value = 1702
$ws.Range("M113").Value = -1080
$ws.Range("N113").Value = -6042

# CRP row 139: 'Weaving a Path' / 'Acacia Spinning Wheel'
$ws.Range("H139").Value = 25100.385
$ws.Range("J139").Value = 25100.385
$ws.Range("L139").Value = 25100.385
$ws.Range("N139").Value = -35380.38499999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 34: 'Fever Pitch' / 'Chamomile Tea'
$ws.Range("H34").Value = 66667536
$ws.Range("I34").Value = 93.5
$ws.Range("J34").Value = 76924060
$ws.Range("K34").Value = 280.5
$ws.Range("L34").Value = 230772180
$ws.Range("M34").Value = -196.5
$ws.Range("N34").Value = -230772348

# CUL row 113: 'Can''t Eat Just One' / 'Night Vinegar'
$ws.Range("H113").Value = 752
$ws.Range("I113").Value = 790
$ws.Range("J113").Value = 741.1429000000001
$ws.Range("K113").Value = 2370
$ws.Range("L113").Value = 2223.4287
$ws.Range("M113").Value = -200
$ws.Range("N113").Value = -6563.4287

# CUL row 132: 'More Mezcal' / 'Cooking Mezcal'
$ws.Range("H132").Value = 52632650
$ws.Range("I132").Value = 62500880
$ws.Range("K132").Value = 562507920
$ws.Range("M132").Value = -562505390

# CUL row 137: 'Creative Chocolate' / 'Gateau au Chocolat'
$ws.Range("H137").Value = 63101.883
$ws.Range("I137").Value = 3180
$ws.Range("J137").Value = 84673.75999999999
$ws.Range("K137").Value = 9540
$ws.Range("L137").Value = 254021.28
$ws.Range("M137").Value = -4440
$ws.Range("N137").Value = -264221.28

# CUL row 138: 'Bring Me Your Tacos' / 'Tacos Al Pastor'
$ws.Range("H138").Value = 1796.3636
$ws.Range("I138").Value = 1796.3636
$ws.Range("K138").Value = 5389.0908
$ws.Range("M138").Value = -249.0907999999999

$ws = $wb.Worksheets.Item("GSM")
# GSM row 132: 'On Board for Lar' / 'Lar Ingot'
$ws.Range("H132").Value = 1611.1923
$ws.Range("I132").Value = 1169.2106
$ws.Range("K132").Value = 3507.6318
$ws.Range("M132").Value = -977.6318000000001

$ws = $wb.Worksheets.Item("LTW")
# LTW row 61: 'Spelling Me Softly' / 'Raptor Leather'
$ws.Range("H61").Value = 1689.2
$ws.Range("I61").Value = 1236.5
$ws.Range("K61").Value = 1236.5
$ws.Range("M61").Value = -1034.5

# LTW row 113: 'Peace in Rest' / 'Atrociraptor Leather'
$ws.Range("H113").Value = 1689.2
$ws.Range("I113").Value = 1236.5
$ws.Range("K113").Value = 1236.5
$ws.Range("M113").Value = 933.5

# LTW row 122: 'Hell on Leather' / 'Gaja Leather'
$ws.Range("H122").Value = 3506.3333
$ws.Range("I122").Value = 3257.6
$ws.Range("J122").Value = 4750
$ws.Range("K122").Value = 9772.799999999999
$ws.Range("L122").Value = 14250
$ws.Range("M122").Value = -7322.799999999999
$ws.Range("N122").Value = -19150

# LTW row 132: 'Tenets of Tanning' / 'Silver Lobo Leather'
$ws.Range("H132").Value = 9469.375
$ws.Range("I132").Value = 18352
$ws.Range("J132").Value = 4139.8
$ws.Range("K132").Value = 55056
$ws.Range("L132").Value = 12419.4
$ws.Range("M132").Value = -52526
$ws.Range("N132").Value = -17479.4

# LTW row 136: 'Respect for Br''aax' / 'Br''aax Leather'
$ws.Range("H136").Value = 2569.5454
$ws.Range("I136").Value = 1464
$ws.Range("J136").Value = 3490.8333
$ws.Range("K136").Value = 4392
$ws.Range("L136").Value = 10472.4999
$ws.Range("M136").Value = -1842
$ws.Range("N136").Value = -15572.4999

# LTW row 138: 'Freezing Toes' / 'Gomphotherium Boots of Striking'
$ws.Range("H138").Value = 50607.25
$ws.Range("I138").Value = 25000
$ws.Range("J138").Value = 54265.43
$ws.Range("K138").Value = 25000
$ws.Range("L138").Value = 54265.43
$ws.Range("M138").Value = -19860
$ws.Range("N138").Value = -64545.43

$ws = $wb.Worksheets.Item("WVR")
# WVR row 113: 'A Tender Table' / 'Pixie Floss'
$ws.Range("H113").Value = 461.03333
$ws.Range("I113").Value = 428.45
$ws.Range("K113").Value = 1285.35
$ws.Range("M113").Value = 884.6500000000001

# WVR row 132: 'Comfy Cabins' / 'Snow Cotton Cloth'
$ws.Range("H132").Value = 6749.5835
$ws.Range("I132").Value = 12217.6
$ws.Range("J132").Value = 2843.8572
$ws.Range("K132").Value = 36652.8
$ws.Range("L132").Value = 8531.571599999999
$ws.Range("M132").Value = -34122.8
$ws.Range("N132").Value = -13591.5716
